# Implements: Added __footings__attribute_map__ as attribute under model to be used.
# This prefixes parameter names with "parameter." and asset/return names with "asset."
# across the step sheets (_step_1, _step_2, _step_3) and the hidden __footings__
# metadata sheet. Also widens the "Uses/Impacts" column on the step sheets to fit
# the longer text.

$wb = $excel.ActiveWorkbook

# --- _step_1 sheet ---------------------------------------------------------
$step1 = $wb.Worksheets.Item("_step_1")
$step1.Range("C8").Value = "[parameter.a, parameter.b]"
$step1.Range("C10").Value = "[asset.ret_1]"
$step1.Range("C14").Value = "asset.ret_1"
$step1.Columns.Item(3).ColumnWidth = 29 - 5/6

# --- _step_2 sheet ---------------------------------------------------------
$step2 = $wb.Worksheets.Item("_step_2")
$step2.Range("C8").Value = "[parameter.c, parameter.d]"
$step2.Range("C10").Value = "[asset.ret_2]"
$step2.Range("C14").Value = "asset.ret_2"
$step2.Columns.Item(3).ColumnWidth = 29 - 5/6

# --- _step_3 sheet ---------------------------------------------------------
$step3 = $wb.Worksheets.Item("_step_3")
$step3.Range("C8").Value = "[asset.ret_1, asset.ret_2]"
$step3.Range("C10").Value = "[asset.ret_3]"
$step3.Range("C14").Value = "asset.ret_3"

# --- hidden __footings__ metadata sheet -------------------------------------
$footings = $wb.Worksheets.Item("__footings__")
$footings.Range("C57").Value = "[asset.ret_1]"
$footings.Range("C58").Value = "[asset.ret_1]"
$footings.Range("C72").Value = "[asset.ret_2]"
$footings.Range("C73").Value = "[asset.ret_2]"
$footings.Range("C87").Value = "[asset.ret_3]"
$footings.Range("C88").Value = "[asset.ret_3]"
